$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on every cell we touch so that numeric-looking
# strings (e.g. "1.000", "0.07645") are preserved verbatim as text, matching
# the original inlineStr cell content instead of being auto-coerced to a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "22.301.06"
$ws.Range("E2").Value = "  -5.16%  "
$ws.Range("D3").Value = "1.565.16"
$ws.Range("E3").Value = "  -5.23%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "289.37"
$ws.Range("E6").Value = "  -3.59%  "
$ws.Range("D7").Value = "0.3746"
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("D8").Value = "49.31"
$ws.Range("E8").Value = "  -2.57%  "
$ws.Range("D9").Value = "0.3413"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").Value = "1.167"
$ws.Range("E10").Value = "  -4.92%  "
$ws.Range("D11").Value = "0.07645"
$ws.Range("E11").Value = "  -5.20%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "21.42"
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("D14").Value = "6.014"
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("D15").Value = "6.936"
$ws.Range("E15").Value = "  -4.67%  "
$ws.Range("D16").Value = "1.561.56"
$ws.Range("E16").Value = "  -5.36%  "
$ws.Range("D17").Value = "0.00001130"
$ws.Range("E17").Value = "  -6.90%  "
$ws.Range("D18").Value = "89.97"
$ws.Range("E18").Value = "  -5.69%  "
$ws.Range("D19").Value = "0.06713"
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "6.237"
$ws.Range("E21").Value = "  -6.09%  "
$ws.Range("D22").Value = "16.56"
$ws.Range("E22").Value = "  -5.23%  "
$ws.Range("D23").Value = "0.5275"
$ws.Range("E23").Value = "  -7.84%  "
$ws.Range("D24").Value = "11.94"
$ws.Range("D25").Value = "22.295.16"
$ws.Range("E25").Value = "  -5.18%  "
$ws.Range("D26").Value = "2.390"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").Value = "2.778"
$ws.Range("E27").Value = "  -8.17%  "
$ws.Range("D28").Value = "20.18"
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("D30").Value = "4.964"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("D31").Value = "125.42"
$ws.Range("E31").Value = "  -4.90%  "
$ws.Range("D32").Value = "1.734.48"
$ws.Range("E32").Value = "  -5.57%  "
$ws.Range("E33").Value = "  +3.36%  "
$ws.Range("D34").Value = "6.201"
$ws.Range("E34").Value = "  -10.30%  "
$ws.Range("D35").Value = "2.007"
$ws.Range("E35").Value = "  -6.14%  "
$ws.Range("D36").Value = "10.06"
$ws.Range("E36").Value = "  -10.61%  "
$ws.Range("D37").Value = "0.08531"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("D38").Value = "0.02534"
$ws.Range("E38").Value = "  -6.98%  "
$ws.Range("D39").Value = "0.2320"
$ws.Range("E39").Value = "  -4.41%  "
$ws.Range("D40").Value = "5.525"
$ws.Range("E40").Value = "  -7.16%  "
$ws.Range("D41").Value = "1.324"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").Value = "0.06408"
$ws.Range("E42").Value = "  -6.25%  "
$ws.Range("D43").Value = "11.70"
$ws.Range("E43").Value = "  -9.56%  "
$ws.Range("D44").Value = "0.6368"
$ws.Range("E44").Value = "  -7.89%  "
$ws.Range("D45").Value = "14.10"
$ws.Range("E45").Value = "  -10.10%  "
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "0.5981"
$ws.Range("E47").Value = "  -6.61%  "
$ws.Range("D48").Value = "3.752"
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("D49").Value = "2.094"
$ws.Range("E49").Value = "  -7.07%  "
$ws.Range("D50").Value = "1.269"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").Value = "123.99"
$ws.Range("E51").Value = "  -2.58%  "
